$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Add the new "SL_ratio" worksheet after "Fit leaf_L"
# ---------------------------------------------------------------------------
$lastIdx = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($lastIdx)
$ws = $wb.Worksheets.Add([System.Type]::Missing, $lastSheet)
$ws.Name = "SL_ratio"

# ---------------------------------------------------------------------------
# 2. Header lines (A1, A2) - python-like docstring
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "def calculate_SL_ratio(leaf_rank):"
$ws.Range("A1").Font.Name = "Courier New"
$ws.Range("A1").Font.Size = 10
$ws.Range("A1").Font.Bold = $true
$ws.Range("A1").Font.Color = 250
$ws.Range("A1").VerticalAlignment = -4108

$a1 = $ws.Range("A1")
$c = $a1.Characters(1, 4)
$c.Font.Name = "Calibri"
$c.Font.Size = 11
$c.Font.Bold = $false
$c.Font.Color = -4105
$c2 = $a1.Characters(5, 19)
$c2.Font.Name = "Courier New"
$c2.Font.Size = 10
$c2.Font.Color = 0
$c2.Font.Bold = $false
$c3 = $a1.Characters(24, 1)
$c3.Font.Name = "Courier New"
$c3.Font.Size = 10
$c3.Font.Color = 6684672
$c3.Font.Bold = $false
$c4 = $a1.Characters(25, 9)
$c4.Font.Name = "Courier New"
$c4.Font.Size = 10
$c4.Font.Color = 6684672
$c4.Font.Italic = $true
$c4.Font.Bold = $false
$c5 = $a1.Characters(34, 2)
$c5.Font.Name = "Courier New"
$c5.Font.Size = 10
$c5.Font.Color = 6684672
$c5.Font.Bold = $false
$c6 = $a1.Characters(35, 1)
$c6.Font.Name = "Courier New"
$c6.Font.Size = 10
$c6.Font.Color = 3355443
$c6.Font.Bold = $false

$ws.Range("A2").Value = '    """ Sheath:Lamina final length ratio according to the rank. Parameters from Dornbush (2011).'
$ws.Range("A2").Font.Name = "Courier New"
$ws.Range("A2").Font.Size = 10
$ws.Range("A2").Font.Color = 0x333333
$ws.Range("A2").Interior.Color = 0xFFFFFF
$ws.Range("A2").VerticalAlignment = -4108

Write-Output "ok step 2"
